$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = 0.2727157963049881
$ws.Range("E2").Value = 17.31518318016685
$ws.Range("F2").Value = 77.46345199258539
